# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts, as published by the gh-pages data generator.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    4  = 636
    5  = 192
    7  = 9635
    8  = 867
    9  = 330
    10 = 1220
    11 = 2254
    12 = 161
    14 = 10
    16 = 274
    17 = 473
    18 = 99
    19 = 258
    20 = 1354
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "全部类型" (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    5  = 636
    6  = 192
    8  = 9635
    9  = 867
    10 = 330
    11 = 1220
    12 = 2254
    13 = 161
    15 = 10
    17 = 274
    18 = 473
    19 = 99
    20 = 258
    21 = 1354
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
